# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Cell B11 on the "Rules" sheet is edited from "R40" to the text "1"
# (its existing cell style / number format, which is General, is left
# untouched - only the stored value changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")

# A plain `Value = "1"` assignment would be auto-coerced by Excel's
# General-format smart typing into the *number* 1. To reproduce a literal
# text "1" (as the original cell held literal text "R40") without
# disturbing the cell's formatting/style, compute the text via a formula
# and then collapse it down to a static value with a values-only paste -
# the same end result a user gets from Paste Special > Values.
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
